$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 7 already holds the constant "mobile / taluko / gaam" pattern
# (A=9510851351, D=1, E=221) that the two new rows reuse. Copy it down
# so the numeric-looking text stays text (same shared-string refs, no
# new cell style), then overwrite only the name/date columns.
$ws.Range("A7:E7").Copy()
$ws.Range("A12").PasteSpecial()
$ws.Range("B12").Value = "સ્વ. ગોવિંદભાઈ ભગવાનભાઇ કાનપરિયા "
$ws.Range("C12").Value = "તારીખ -૦૭-૦૭-૨૦૨૫ વાર સોમવાર "

$ws.Range("A7:E7").Copy()
$ws.Range("A13").PasteSpecial()
$ws.Range("B13").Value = "સ્વ. રાજાભાઈ અમરાભાઇ કણબી "
$ws.Range("C13").Value = "તારીખ -૧૨-૦૭-૨૦૨૫ વાર - શનિવાર"

$excel.CutCopyMode = $false
